$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Дефект1")

# Set the new value in B9 (this also adds the shared string and selects the cell in Excel UI)
$ws.Range("B9").Value = "Just checking it out. Tanya"

# Make sure the sheet is active and the selection matches B9
$ws.Activate()
$ws.Range("B9").Select()
